$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.518329262733459
$ws.Range("B1").Value = 1.843080997467041
$ws.Range("C1").Value = 1.961581945419312
$ws.Range("D1").Value = 2.265778779983521
$ws.Range("E1").Value = 2.852095365524292
